$wb = $excel.ActiveWorkbook

# Sheet "展览" (first sheet)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F8").Value = 59
$ws1.Range("F9").Value = 8382
$ws1.Range("F10").Value = 780
$ws1.Range("F11").Value = 308
$ws1.Range("F13").Value = 889
$ws1.Range("F14").Value = 73
$ws1.Range("F17").Value = 154
$ws1.Range("F20").Value = 931

# Sheet "全部类型" (fourth sheet)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F10").Value = 59
$ws4.Range("F11").Value = 8383
$ws4.Range("F12").Value = 780
$ws4.Range("F13").Value = 308
$ws4.Range("F15").Value = 889
$ws4.Range("F16").Value = 73
$ws4.Range("F19").Value = 154
$ws4.Range("F22").Value = 931
